$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 5277.1665
$ws.Range("I82").Value = 1903.7142
$ws.Range("K82").Value = 5711.142599999999
$ws.Range("M82").Value = -5305.142599999999
$ws.Range("H85").Value = 5277.1665
$ws.Range("I85").Value = 1903.7142
$ws.Range("K85").Value = 5711.142599999999
$ws.Range("M85").Value = -4307.142599999999
$ws.Range("H101").Value = 700
$ws.Range("I101").Value = 700
$ws.Range("K101").Value = 2100
$ws.Range("M101").Value = -478
$ws.Range("H116").Value = 2200557.8
$ws.Range("I116").Value = 3796862.8
$ws.Range("J116").Value = 5638.5
$ws.Range("K116").Value = 3796862.8
$ws.Range("L116").Value = 5638.5
$ws.Range("M116").Value = -3793420.8
$ws.Range("N116").Value = -12522.5
$ws.Range("H129").Value = 1232.1
$ws.Range("I129").Value = 553.125
$ws.Range("K129").Value = 1659.375
$ws.Range("M129").Value = 3340.625
$ws.Range("H132").Value = 151628.61
$ws.Range("I132").Value = 336201.97
$ws.Range("J132").Value = 13198.583
$ws.Range("K132").Value = 1008605.91
$ws.Range("L132").Value = 39595.749
$ws.Range("M132").Value = -1006075.91
$ws.Range("N132").Value = -44655.749
$ws.Range("H135").Value = 7199.375
$ws.Range("I135").Value = 2739
$ws.Range("K135").Value = 24651
$ws.Range("M135").Value = -22116
$ws.Range("H137").Value = 4786.75
$ws.Range("I137").Value = 1633.3334
$ws.Range("K137").Value = 4900.0002
$ws.Range("M137").Value = -2350.0002
$ws.Range("H138").Value = 6820.7896
$ws.Range("J138").Value = 8537.885
$ws.Range("L138").Value = 25613.655
$ws.Range("N138").Value = -35893.655
$ws.Range("H140").Value = 73710.48
$ws.Range("J140").Value = 74357
$ws.Range("L140").Value = 74357
$ws.Range("N140").Value = -84717

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12912.046
$ws.Range("I61").Value = 12653.143
$ws.Range("K61").Value = 12653.143
$ws.Range("M61").Value = -12441.143
$ws.Range("H74").Value = 881.7
$ws.Range("I74").Value = 574
$ws.Range("J74").Value = 1599.6666
$ws.Range("K74").Value = 574
$ws.Range("L74").Value = 1599.6666
$ws.Range("M74").Value = 300
$ws.Range("N74").Value = -3347.6666
$ws.Range("H77").Value = 881.7
$ws.Range("I77").Value = 574
$ws.Range("J77").Value = 1599.6666
$ws.Range("K77").Value = 2870
$ws.Range("L77").Value = 7998.333000000001
$ws.Range("M77").Value = 1498
$ws.Range("N77").Value = -16734.333
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 5399.9614
$ws.Range("I132").Value = 2280.611
$ws.Range("J132").Value = 12418.5
$ws.Range("K132").Value = 6841.833
$ws.Range("L132").Value = 37255.5
$ws.Range("M132").Value = -4311.833
$ws.Range("N132").Value = -42315.5
$ws.Range("H136").Value = 12912.046
$ws.Range("I136").Value = 12653.143
$ws.Range("K136").Value = 37959.429
$ws.Range("M136").Value = -35409.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3683.1155
$ws.Range("I86").Value = 3600.3845
$ws.Range("K86").Value = 3600.3845
$ws.Range("M86").Value = -2477.3845
$ws.Range("H89").Value = 3683.1155
$ws.Range("I89").Value = 3600.3845
$ws.Range("K89").Value = 18001.9225
$ws.Range("M89").Value = -12385.9225
$ws.Range("H105").Value = 2609.8948
$ws.Range("I105").Value = 3040.6667
$ws.Range("K105").Value = 3040.6667
$ws.Range("M105").Value = -1293.6667
$ws.Range("H134").Value = 3885.8823
$ws.Range("I134").Value = 3296.8333
$ws.Range("K134").Value = 9890.499899999999
$ws.Range("M134").Value = -7355.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1531.2667
$ws.Range("I16").Value = 1447.8334
$ws.Range("K16").Value = 1447.8334
$ws.Range("M16").Value = -1160.8334
$ws.Range("H74").Value = 79730.164
$ws.Range("J74").Value = 79730.164
$ws.Range("L74").Value = 79730.164
$ws.Range("N74").Value = -81478.164
$ws.Range("H77").Value = 79730.164
$ws.Range("J77").Value = 79730.164
$ws.Range("L77").Value = 239190.492
$ws.Range("N77").Value = -247926.492
$ws.Range("H95").Value = 47897.4
$ws.Range("J95").Value = 47897.4
$ws.Range("L95").Value = 47897.4
$ws.Range("N95").Value = -53389.4
$ws.Range("H105").Value = 2067187.8
$ws.Range("I105").Value = 2841608
$ws.Range("K105").Value = 2841608
$ws.Range("M105").Value = -2839861
$ws.Range("H107").Value = 729754.25
$ws.Range("I107").Value = 1213095.8
$ws.Range("J107").Value = 4741.9
$ws.Range("K107").Value = 1213095.8
$ws.Range("L107").Value = 4741.9
$ws.Range("M107").Value = -1211175.8
$ws.Range("N107").Value = -8581.9
$ws.Range("H113").Value = 1531.2667
$ws.Range("I113").Value = 1447.8334
$ws.Range("K113").Value = 1447.8334
$ws.Range("M113").Value = 722.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 973849.5
$ws.Range("I5").Value = 2336
$ws.Range("J5").Value = 1556757.6
$ws.Range("K5").Value = 7008
$ws.Range("L5").Value = 4670272.800000001
$ws.Range("M5").Value = -6896
$ws.Range("N5").Value = -4670496.800000001
$ws.Range("H68").Value = 1297479.5
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1945719.2
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 5837157.6
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -5838779.6
$ws.Range("H71").Value = 1297479.5
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1945719.2
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 17511472.8
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -17519584.8
$ws.Range("H93").Value = 2504.5
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 2805.4
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 8416.200000000001
$ws.Range("M93").Value = -1128
$ws.Range("N93").Value = -12160.2
$ws.Range("H107").Value = 1526.7142
$ws.Range("J107").Value = 258.75
$ws.Range("L107").Value = 776.25
$ws.Range("N107").Value = -4616.25
$ws.Range("H109").Value = 16234.5
$ws.Range("I109").Value = 1632
$ws.Range("J109").Value = 24996
$ws.Range("K109").Value = 4896
$ws.Range("L109").Value = 74988
$ws.Range("M109").Value = -3856
$ws.Range("N109").Value = -77068
$ws.Range("H127").Value = 1416
$ws.Range("J127").Value = 1416
$ws.Range("L127").Value = 4248
$ws.Range("N127").Value = -14168
$ws.Range("H132").Value = 5333.3335
$ws.Range("J132").Value = 5333.3335
$ws.Range("L132").Value = 48000.0015
$ws.Range("N132").Value = -53060.0015
$ws.Range("H135").Value = 973849.5
$ws.Range("I135").Value = 2336
$ws.Range("J135").Value = 1556757.6
$ws.Range("K135").Value = 21024
$ws.Range("L135").Value = 14010818.4
$ws.Range("M135").Value = -18489
$ws.Range("N135").Value = -14015888.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 59915
$ws.Range("I138").Value = 59948
$ws.Range("J138").Value = 59849
$ws.Range("K138").Value = 59948
$ws.Range("L138").Value = 59849
$ws.Range("M138").Value = -54808
$ws.Range("N138").Value = -70129

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3286.5454
$ws.Range("I40").Value = 2072.7856
$ws.Range("K40").Value = 2072.7856
$ws.Range("M40").Value = -1936.7856
$ws.Range("H58").Value = 17046.5
$ws.Range("J58").Value = 30000
$ws.Range("L58").Value = 30000
$ws.Range("N58").Value = -30520
$ws.Range("H132").Value = 5783.04
$ws.Range("I132").Value = 3515.6
$ws.Range("J132").Value = 7294.6665
$ws.Range("K132").Value = 10546.8
$ws.Range("L132").Value = 21883.9995
$ws.Range("M132").Value = -8016.799999999999
$ws.Range("N132").Value = -26943.9995
$ws.Range("H136").Value = 10000
$ws.Range("I136").Value = 10000
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -27450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41680890
$ws.Range("I132").Value = 2806.2856
$ws.Range("J132").Value = 100030200
$ws.Range("K132").Value = 8418.856800000001
$ws.Range("L132").Value = 300090600
$ws.Range("M132").Value = -5888.856800000001
$ws.Range("N132").Value = -300095660
$ws.Range("H136").Value = 5117.24
$ws.Range("I136").Value = 4133
$ws.Range("K136").Value = 12399
$ws.Range("M136").Value = -9849

Write-Host "Applied all profit updates"